# The edited sheet is "Valve_8.0_600_4" (already the active sheet/tab in
# the source workbook - tabSelected="1", activeTab="10").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valve_8.0_600_4")
$ws.Activate()

# Update the B3:B12 measurement values (column B) per the diff.
$ws.Range("B3").Value = 100
$ws.Range("B4").Value = 200
$ws.Range("B5").Value = 300
$ws.Range("B6").Value = 400
$ws.Range("B7").Value = 500
$ws.Range("B8").Value = 600
$ws.Range("B9").Value = 700
$ws.Range("B10").Value = 800
$ws.Range("B11").Value = 900
$ws.Range("B12").Value = 1000

# Move the active selection to I20, matching the saved cursor position.
[void]$ws.Range("I20").Select()
